$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking Price values (column D stores prices as text)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated values from the latest crypto data pull
$ws.Range("D2").Value = "65.129.06"
$ws.Range("E2").Value = "  +3.92%  "
$ws.Range("D3").Value = "2.561.23"
$ws.Range("E3").Value = "  +4.04%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "582.32"
$ws.Range("E5").Value = "  +1.62%  "
$ws.Range("D6").Value = "153.58"
$ws.Range("E6").Value = "  +4.25%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +1.46%  "
$ws.Range("D9").Value = "2.566.76"
$ws.Range("E9").Value = "  +4.18%  "
$ws.Range("E10").Value = "  +1.56%  "
$ws.Range("E11").Value = "  -1.66%  "
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("E13").Value = "  +1.04%  "
$ws.Range("D14").Value = "29.36"
$ws.Range("E14").Value = "  +1.37%  "
$ws.Range("E15").Value = "  +2.72%  "
$ws.Range("D16").Value = "3.022.34"
$ws.Range("E16").Value = "  +3.75%  "
$ws.Range("D17").Value = "64.668.77"
$ws.Range("E17").Value = "  +3.25%  "
$ws.Range("D18").Value = "2.557.12"
$ws.Range("E18").Value = "  +3.58%  "
$ws.Range("E19").Value = "  +1.86%  "
$ws.Range("D20").Value = "11.07"
$ws.Range("E20").Value = "  +1.58%  "
$ws.Range("D21").Value = "354.21"
$ws.Range("E21").Value = "  +8.70%  "
$ws.Range("E22").Value = "  +4.31%  "
$ws.Range("D23").Value = "2.24"
$ws.Range("E23").Value = "  +3.16%  "
$ws.Range("E24").Value = "  +0.00%  "
$ws.Range("D25").Value = "10.12"
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("D26").Value = "66.24"
$ws.Range("E26").Value = "  +1.36%  "
$ws.Range("D27").Value = "636.55"
$ws.Range("E27").Value = "  -0.89%  "
$ws.Range("E28").Value = "  +8.08%  "
$ws.Range("D29").Value = "2.682.11"
$ws.Range("E29").Value = "  +3.72%  "
$ws.Range("E30").Value = "  +4.91%  "
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").Value = "8.15"
$ws.Range("E32").Value = "  +2.74%  "
$ws.Range("E33").Value = "  +2.88%  "
$ws.Range("E34").Value = "  +4.66%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").Value = "1.57"
$ws.Range("E36").Value = "  +3.73%  "
$ws.Range("E37").Value = "  +3.38%  "
$ws.Range("D38").Value = "5.62"
$ws.Range("E38").Value = "  +6.06%  "
$ws.Range("D39").Value = "2.87"
$ws.Range("E39").Value = "  +5.29%  "
$ws.Range("B40").Value = "EthereumClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D40").Value = "19.18"
$ws.Range("E40").Value = "  +3.37%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "154.69"
$ws.Range("E41").Value = "  +1.95%  "
$ws.Range("E42").Value = "  +1.50%  "
$ws.Range("E43").Value = "  +5.29%  "
$ws.Range("D44").Value = "42.14"
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("D45").Value = "161.45"
$ws.Range("E45").Value = "  +5.49%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0307"
$ws.Range("E46").Value = "  -1.91%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").Value = "  -0.03%  "
$ws.Range("D48").Value = "15.71"
$ws.Range("E48").Value = "  +2.67%  "
$ws.Range("E49").Value = "  +3.39%  "
$ws.Range("D50").Value = "21.71"
$ws.Range("D51").Value = "0.632"
$ws.Range("E51").Value = "  +4.28%  "
